$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Games")

$startRow = 54
$startId = 53

for ($i = 0; $i -lt 9; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $startId + $i
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 3
    $ws.Cells.Item($r, 5).Value = "14.01.2020"
}
